$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: only the count (B) changes
$ws.Range("B24").Value = 57612

# Row 25: becomes the old row 27 species data, with new Id/Taxonsorteringsordning
$ws.Range("A25").Value = 112194813
$ws.Range("B25").Value = 90087
$ws.Range("E25").Value = 3884
$ws.Range("F25").Value = "Hasselticka"
$ws.Range("G25").Value = "Dichomitus campestris"
$ws.Range("H25").Value = "(Quél.) Domański & Orlicz"

# Row 26: only the count (B) changes
$ws.Range("B26").Value = 56446

# Row 27: becomes the old row 28 species data, with updated Id/Taxonsorteringsordning and times
$ws.Range("A27").Value = 112196967
$ws.Range("B27").Value = 43473
$ws.Range("E27").Value = 101735
$ws.Range("F27").Value = "Jättesvampmal"
$ws.Range("G27").Value = "Scardia boletella"
$ws.Range("H27").Value = "(Fabricius, 1794)"
$ws.Range("Z27").Value = "12:30"
$ws.Range("AB27").Value = "12:30"

# Row 28: becomes the old row 25 species data, with the old row 25 Id/Taxonsorteringsordning and times
$ws.Range("A28").Value = 112195278
$ws.Range("B28").Value = 8377
$ws.Range("E28").Value = 106545
$ws.Range("F28").Value = "Mindre märgborre"
$ws.Range("G28").Value = "Tomicus minor"
$ws.Range("H28").Value = "(Hartig, 1834)"
$ws.Range("Z28").Value = "11:00"
$ws.Range("AB28").Value = "11:00"

# Row 29: only the count (B) changes
$ws.Range("B29").Value = 60196
